# Apply the changes described in the diff:
# - Metadata!B8: updated Date string
# - Concepts!B2:B7: renamed codes

$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2021-11-30T21:05:25+00:00"

$wsConcepts = $wb.Worksheets.Item("Concepts")
$wsConcepts.Range("B2").Value = "MMG"
$wsConcepts.Range("B3").Value = "DYSTM"
$wsConcepts.Range("B4").Value = "RHAB"
$wsConcepts.Range("B5").Value = "MYOPC"
$wsConcepts.Range("B6").Value = "MYASC"
$wsConcepts.Range("B7").Value = "HYPM"
